$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column B (Buying Opportunity) values for rows 2-8
$ws.Range("B2").Value = "NSE:DPABHUSHAN"
$ws.Range("B3").Value = "NSE:GILT5YBEES"
$ws.Range("B4").Value = "NSE:MANINDS"
$ws.Range("B5").Value = "NSE:NACLIND"
$ws.Range("B6").Value = "NSE:OSWALAGRO"
$ws.Range("B7").Value = "NSE:PRECOT"
$ws.Range("B8").Value = "NSE:RKDL"

# Update column C (support Zone) values for rows 2-4, clear rows 5-18
$ws.Range("C2").Value = "NSE:AFFLE"
$ws.Range("C3").Value = "NSE:ANUP"
$ws.Range("C4").Value = "NSE:M&M"
$ws.Range("C5:C18").ClearContents()

# Clear column E (Short buildup) for rows 2-4
$ws.Range("E2:E4").ClearContents()

# Delete rows 9-18 entirely (shift cells up)
$ws.Rows("9:18").Delete()
